$wb = $excel.ActiveWorkbook

# Duplicate the "2021-Q4" sheet (same column layout/styling as the new
# quarter) and place the copy immediately before "总计", then rename it.
$src = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Update the new sheet's two data rows with the 2022-Q1 figures. These
# columns hold number-looking text (matching the existing "1.04"/"2.29"
# style cells already on the sheet), so force the cells to Text before
# writing, then drop back to the Normal style so no stray number format
# is left behind on the cell.
foreach ($r in 2,3) {
    foreach ($col in "D","E","F","G") {
        $cell = $newSheet.Range($col + $r)
        $cell.NumberFormat = "@"
    }
}
$newSheet.Range("D2").Value = "1.84"
$newSheet.Range("E2").Value = "86.48"
$newSheet.Range("F2").Value = "4.17"
$newSheet.Range("G2").Value = "0.0767"
$newSheet.Range("H2").Value = 3

$newSheet.Range("D3").Value = "1.84"
$newSheet.Range("E3").Value = "86.48"
$newSheet.Range("F3").Value = "4.17"
$newSheet.Range("G3").Value = "0.0767"
$newSheet.Range("H3").Value = 3

foreach ($r in 2,3) {
    foreach ($col in "D","E","F","G") {
        $newSheet.Range($col + $r).Style = "Normal"
    }
}

# Insert a new top data row in "总计" for 2022-Q1, pushing the existing
# rows down, then fix up the leading index column. (Re-fetch the sheet --
# the earlier $total reference became stale once the new sheet was
# inserted in front of it.)
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The insert copies the header row's bold/bordered formatting onto the
# new row; reset it back to the plain look used by the other data rows.
$total.Range("B2:D2").Style = "Normal"
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.15

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
